# Generate Report for Handoff
# Adds two new localization entries (23d785a5-... and df0478a4-...) ahead of
# the existing ".localization-config" row on all three worksheets
# (Overview, zh-cn, de-de), shifting ".localization-config" down by two rows.

$wb = $excel.ActiveWorkbook

$mdCommit = "8fac9b1e550a535ad10baa7bca75bff31cb923ea"
$zhCommit = "50058f9b0fcab5b8799280224081f6f55d21ae2a"
$deCommit = "8414ef9e5e2641e78cc0e7556e30999ad94fa270"

$linkColor = 15570276   # BGR encoding of RGB(0x64,0x95,0xED) == style "HyperLink" font color
$linkUnderline = 2      # xlUnderlineStyleSingle

function Set-LinkLook($range) {
    $range.Font.Underline = $linkUnderline
    $range.Font.Color = $linkColor
}

function Update-Hyperlink($ws, $cellAddr, $newUrl, $newDisplay) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellAddr) {
            $h.Address = $newUrl
            $h.TextToDisplay = $newDisplay
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4 now becomes the 23d785a5 entry (re-use the existing hyperlink/rId)
$wsOverview.Range("A4").Value = "23d785a5-46aa-4652-9eb4-47fb758425bd.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
Update-Hyperlink $wsOverview '$A$4' ("https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/23d785a5-46aa-4652-9eb4-47fb758425bd.md") "23d785a5-46aa-4652-9eb4-47fb758425bd.md"

# Row 5 - new df0478a4 entry
$wsOverview.Range("A5").Value = "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md", "", "", "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md") | Out-Null
Set-LinkLook $wsOverview.Range("A5")

# Row 6 - shifted ".localization-config" entry
$wsOverview.Range("A6").Value = ".localization-config"
$wsOverview.Range("B6").Value = "Not to be localized"
$wsOverview.Range("C6").Value = "Not to be localized"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/.localization-config", "", "", ".localization-config") | Out-Null
Set-LinkLook $wsOverview.Range("A6")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 4 now becomes the 23d785a5 entry
$wsZh.Range("A4").Value = "23d785a5-46aa-4652-9eb4-47fb758425bd.md"
$wsZh.Range("B4").Value = "Ready for handoff"
Update-Hyperlink $wsZh '$A$4' ("https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/23d785a5-46aa-4652-9eb4-47fb758425bd.md") "23d785a5-46aa-4652-9eb4-47fb758425bd.md"

$wsZh.Range("C4").Value = "23d785a5-46aa-4652-9eb4-47fb758425bd.cbb0ee4354197207cc1fe370a1a8ff2300d647ec.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/23d785a5-46aa-4652-9eb4-47fb758425bd.cbb0ee4354197207cc1fe370a1a8ff2300d647ec.zh-cn.xlf", "", "", "23d785a5-46aa-4652-9eb4-47fb758425bd.cbb0ee4354197207cc1fe370a1a8ff2300d647ec.zh-cn.xlf") | Out-Null
Set-LinkLook $wsZh.Range("C4")
$wsZh.Range("D4").Value = "2016-01-27 07:33:09"
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"

# Row 5 - new df0478a4 entry
$wsZh.Range("A5").Value = "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md"
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md", "", "", "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md") | Out-Null
Set-LinkLook $wsZh.Range("A5")

$wsZh.Range("C5").Value = "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.713b575872537f9075ef300a365bb1924d255190.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.713b575872537f9075ef300a365bb1924d255190.zh-cn.xlf", "", "", "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.713b575872537f9075ef300a365bb1924d255190.zh-cn.xlf") | Out-Null
Set-LinkLook $wsZh.Range("C5")
$wsZh.Range("D5").Value = "2016-01-27 07:33:09"
$wsZh.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"

# Row 6 - shifted ".localization-config" entry
$wsZh.Range("A6").Value = ".localization-config"
$wsZh.Range("B6").Value = "Not to be localized"
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/.localization-config", "", "", ".localization-config") | Out-Null
Set-LinkLook $wsZh.Range("A6")
$wsZh.Range("D6").Value = "0001-01-01 00:00:00"
$wsZh.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G6").Value = "0001-01-01 00:00:00"
$wsZh.Range("H6").Value = "Ignored"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 4 now becomes the 23d785a5 entry
$wsDe.Range("A4").Value = "23d785a5-46aa-4652-9eb4-47fb758425bd.md"
$wsDe.Range("B4").Value = "Ready for handoff"
Update-Hyperlink $wsDe '$A$4' ("https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/23d785a5-46aa-4652-9eb4-47fb758425bd.md") "23d785a5-46aa-4652-9eb4-47fb758425bd.md"

$wsDe.Range("C4").Value = "23d785a5-46aa-4652-9eb4-47fb758425bd.cbb0ee4354197207cc1fe370a1a8ff2300d647ec.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/23d785a5-46aa-4652-9eb4-47fb758425bd.cbb0ee4354197207cc1fe370a1a8ff2300d647ec.de-de.xlf", "", "", "23d785a5-46aa-4652-9eb4-47fb758425bd.cbb0ee4354197207cc1fe370a1a8ff2300d647ec.de-de.xlf") | Out-Null
Set-LinkLook $wsDe.Range("C4")
$wsDe.Range("D4").Value = "2016-01-27 07:33:19"
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"

# Row 5 - new df0478a4 entry
$wsDe.Range("A5").Value = "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md"
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md", "", "", "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.md") | Out-Null
Set-LinkLook $wsDe.Range("A5")

$wsDe.Range("C5").Value = "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.713b575872537f9075ef300a365bb1924d255190.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.713b575872537f9075ef300a365bb1924d255190.de-de.xlf", "", "", "df0478a4-87d5-45a4-b5fb-b1a6adc0bcf6.713b575872537f9075ef300a365bb1924d255190.de-de.xlf") | Out-Null
Set-LinkLook $wsDe.Range("C5")
$wsDe.Range("D5").Value = "2016-01-27 07:33:19"
$wsDe.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"

# Row 6 - shifted ".localization-config" entry
$wsDe.Range("A6").Value = ".localization-config"
$wsDe.Range("B6").Value = "Not to be localized"
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/.localization-config", "", "", ".localization-config") | Out-Null
Set-LinkLook $wsDe.Range("A6")
$wsDe.Range("D6").Value = "0001-01-01 00:00:00"
$wsDe.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDe.Range("H6").Value = "Ignored"

Write-Host "Report generated for handoff."
